# Update "想去人数" (want-to-go count) figures on the 展览 and 全部类型 sheets
# to reflect the regenerated gh-pages output.

$wb = $excel.ActiveWorkbook

# --- Sheet: 展览 ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F7").Value = 182
$ws1.Range("F9").Value = 6381
$ws1.Range("F12").Value = 123
$ws1.Range("F13").Value = 5569
$ws1.Range("F25").Value = 3992
$ws1.Range("F26").Value = 14

# --- Sheet: 全部类型 ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F8").Value = 182
$ws4.Range("F10").Value = 6381
$ws4.Range("F13").Value = 123
$ws4.Range("F14").Value = 5569
$ws4.Range("F26").Value = 3992
$ws4.Range("F28").Value = 14
